# Add a new column to the existing table (table_test), extending it from J to K
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add() | Out-Null

# Header for the new test_xlr_n_percent column
$ws.Range("K2").Value = "test_xlr_n_percent"

# Data rows: "n (percent%)" strings for rows 3..34
$ws.Range("K3").Value = "1 (3%)"
$ws.Range("K4").Value = "2 (6%)"
$ws.Range("K5").Value = "3 (9%)"
$ws.Range("K6").Value = "4 (12%)"
$ws.Range("K7").Value = "5 (16%)"
$ws.Range("K8").Value = "6 (19%)"
$ws.Range("K9").Value = "7 (22%)"
$ws.Range("K10").Value = "8 (25%)"
$ws.Range("K11").Value = "9 (28%)"
$ws.Range("K12").Value = "10 (31%)"
$ws.Range("K13").Value = "11 (34%)"
$ws.Range("K14").Value = "12 (38%)"
$ws.Range("K15").Value = "13 (41%)"
$ws.Range("K16").Value = "14 (44%)"
$ws.Range("K17").Value = "15 (47%)"
$ws.Range("K18").Value = "16 (50%)"
$ws.Range("K19").Value = "17 (53%)"
$ws.Range("K20").Value = "18 (56%)"
$ws.Range("K21").Value = "19 (59%)"
$ws.Range("K22").Value = "20 (62%)"
$ws.Range("K23").Value = "21 (66%)"
$ws.Range("K24").Value = "22 (69%)"
$ws.Range("K25").Value = "23 (72%)"
$ws.Range("K26").Value = "24 (75%)"
$ws.Range("K27").Value = "25 (78%)"
$ws.Range("K28").Value = "26 (81%)"
$ws.Range("K29").Value = "27 (84%)"
$ws.Range("K30").Value = "28 (88%)"
$ws.Range("K31").Value = "29 (91%)"
$ws.Range("K32").Value = "30 (94%)"
$ws.Range("K33").Value = "31 (97%)"
$ws.Range("K34").Value = "32 (100%)"

# Widen column J slightly (matches authored width of 15.1640625 characters)
$ws.Columns.Item(10).ColumnWidth = 14.3

# Update the active selection to K9, matching the authored selection state
$ws.Range("K9").Select() | Out-Null
